# Add a "Format: v0.1.0" label row under the sheet title, matching the
# newest InOutModule template used across the example workbooks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2: B2 = "Format:" (right aligned), C2 = "v0.1.0"
$ws.Range("B2").Value = "Format:"
$ws.Range("C2").Value = "v0.1.0"

# Match the row height used by the rest of the sheet's data rows.
$ws.Rows("2").RowHeight = 18.75

# Style B2: italic 11pt Aptos, right aligned, vertically centered.
$fontB2 = $ws.Range("B2").Font
$fontB2.Italic = $true
$fontB2.Size = 11
$fontB2.Name = "Aptos"
$ws.Range("B2").VerticalAlignment = -4108   # xlCenter
$ws.Range("B2").HorizontalAlignment = -4152 # xlRight

# Style C2: italic 11pt Aptos, vertically centered.
$fontC2 = $ws.Range("C2").Font
$fontC2.Italic = $true
$fontC2.Size = 11
$fontC2.Name = "Aptos"
$ws.Range("C2").VerticalAlignment = -4108   # xlCenter
